$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.851.52"
$ws.Range("E2").Value = "  +7.91%  "

$ws.Range("D3").Value = "1.811.27"
$ws.Range("E3").Value = "  +4.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.54"
$ws.Range("E5").Value = "  +3.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4942"
$ws.Range("E7").Value = "  +2.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2786"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06415"
$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("D10").Value = "1.815.93"
$ws.Range("E10").Value = "  +5.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.78"
$ws.Range("E11").Value = "  +5.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07127"
$ws.Range("E12").Value = "  +3.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6482"
$ws.Range("E13").Value = "  +6.75%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.48"
$ws.Range("E14").Value = "  +9.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.704"
$ws.Range("E15").Value = "  +5.23%  "

$ws.Range("D16").Value = "28.839.30"
$ws.Range("E16").Value = "  +8.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9992"
$ws.Range("E17").Value = "  +0.11%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007391"
$ws.Range("E18").Value = "  +3.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.25"
$ws.Range("E20").Value = "  +7.20%  "

$ws.Range("D21").Value = "2.048.93"
$ws.Range("E21").Value = "  +5.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.597"
$ws.Range("E22").Value = "  +3.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.912"
$ws.Range("E23").Value = "  +3.90%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.360"
$ws.Range("E24").Value = "  +5.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.66"
$ws.Range("E25").Value = "  +4.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "133.73"
$ws.Range("E26").Value = "  +26.09%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.64"
$ws.Range("E27").Value = "  +9.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.889"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.391"
$ws.Range("E29").Value = "  +0.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.167"
$ws.Range("E30").Value = "  +4.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08358"
$ws.Range("E31").Value = "  +5.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.851"
$ws.Range("E32").Value = "  +4.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04946"
$ws.Range("E33").Value = "  +10.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.093"
$ws.Range("E34").Value = "  +8.56%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6774"
$ws.Range("E35").Value = "  +9.47%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.695"
$ws.Range("E36").Value = "  +3.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.292"
$ws.Range("E37").Value = "  +13.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.768"
$ws.Range("E38").Value = "  +13.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9548"
$ws.Range("E39").Value = "  +3.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.051"
$ws.Range("E40").Value = "  +7.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01592"
$ws.Range("E41").Value = "  +6.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9992"
$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.54"
$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4091"
$ws.Range("E44").Value = "  +6.72%  "

$ws.Range("E45").Value = "  +5.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1224"
$ws.Range("E46").Value = "  +5.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05516"
$ws.Range("E47").Value = "  +2.40%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.128"
$ws.Range("E48").Value = "  +3.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.67"
$ws.Range("E49").Value = "  +5.35%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.317"
$ws.Range("E50").Value = "  +7.23%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3629"
$ws.Range("E51").Value = "  +8.01%  "
